# Weekly update: insert a new price-record row for
# "Bruselas (repollito)" at Vega Central Mapocho de Santiago, shifting the
# existing rows 71:98 down to 72:99 and adding a new row 71 with the
# latest week's data (dimension grows from A1:R98 to A1:R99).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one new row above the current row 71 (pushes 71:98 -> 72:99)
$ws.Rows.Item(71).Insert()

# Fill in the new row 71 with the new weekly record
$ws.Range("A71").Value = 9
$ws.Range("B71").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C71").Value = "Metropolitana"
$ws.Range("D71").Value = 45119
$ws.Range("E71").Value = 13
$ws.Range("F71").Value = 100112035
$ws.Range("G71").Value = "Bruselas (repollito)"
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 52
$ws.Range("K71").Value = 18000
$ws.Range("L71").Value = 20000
$ws.Range("M71").Value = 19000
$ws.Range("N71").Value = "$/malla 15 kilos"
$ws.Range("O71").Value = "Provincia de Quillota"
$ws.Range("P71").Value = 1267
$ws.Range("Q71").Value = 15
$ws.Range("R71").Value = "Hortaliza"
